# qa_template_map.xlsx -- add 3 new blank-template fields under the
# "documents" sheet section (effects_data, tk_params, httk_data) that were
# missing, causing duplicate-field bugs downstream. This inserts 3 rows
# right after the existing "documents" block (row 29) and before the
# "series" block (old row 30), pushing everything below down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 30:32 (existing rows 30.. shift down to 33..)
$ws.Rows("30:32").Insert() | Out-Null

# New row 30: documents / effects_data / effects_data
$ws.Range("A30").Value = "documents"
$ws.Range("B30").Value = "effects_data"
$ws.Range("C30").Value = "effects_data"

# New row 31: documents / tk_params / tk_params
$ws.Range("A31").Value = "documents"
$ws.Range("B31").Value = "tk_params"
$ws.Range("C31").Value = "tk_params"

# New row 32: documents / httk_data / httk_data
$ws.Range("A32").Value = "documents"
$ws.Range("B32").Value = "httk_data"
$ws.Range("C32").Value = "httk_data"

# The sheet's hidden AutoFilter-database defined name tracked a specific
# range; keep it consistent with the 3-row shift caused by the insert above.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$C`$107"

# Leave the selection where the edit ended, matching the author's session.
$ws.Range("B32").Select() | Out-Null
